$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.133.36'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '1.573.03'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '207.26'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').Value = '0.490'
$ws.Range('E6').Value = '  -1.58%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').Value = '22.27'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('E10').Value = '  -0.49%  '
$ws.Range('D11').Value = '0.0869'
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').Value = '1.796.55'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').Value = '1.570.82'
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('E15').Value = '  -1.56%  '
$ws.Range('D16').Value = '27.157.79'
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').Value = '62.25'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').Value = '7.37'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').Value = '214.53'
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('E20').Value = '  -1.09%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').Value = '9.46'
$ws.Range('E23').Value = '  -3.68%  '
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('D25').Value = '152.41'
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('D26').Value = '6.70'
$ws.Range('E26').Value = '  -3.46%  '
$ws.Range('D27').Value = '14.93'
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('E30').Value = '  -3.27%  '
$ws.Range('D31').Value = '0.0463'
$ws.Range('E31').Value = '  -2.10%  '
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('D33').Value = '1.399.72'
$ws.Range('E33').Value = '  +1.67%  '
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('D37').Value = '0.942'
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('D38').Value = '0.0164'
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('E39').Value = '  -1.71%  '
$ws.Range('E40').Value = '  -3.64%  '
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('E42').Value = '  +3.99%  '
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('E45').Value = '  +1.08%  '
$ws.Range('D46').Value = '63.74'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('D47').Value = '1.709.30'
$ws.Range('E47').Value = '  -0.69%  '
$ws.Range('D48').Value = '85.51'
$ws.Range('E48').Value = '  +0.19%  '
$ws.Range('D49').Value = '0.0₇0976'
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('D51').Value = '0.0493'
$ws.Range('E51').Value = '  -0.56%  '
